$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new intro text above the table
$ws.Range("F2").Value = "This smart paking system is designed for a pre setup group of people such a hotell where you have users details"

# Mark idea 1 and idea 4 as DONE
$ws.Range("C6").Value = "DONE"
$ws.Range("C9").Value = "DONE"

# Add new idea row 12
$ws.Range("F12").Value = "Admin Login Area To Set Price?"

# Update selection
$ws.Range("Q33").Select()
